$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking strings in D/E/G to remain text (matches source format),
# then strip the temporary text-number-format so cells end up unstyled, like the original.
$deRange = $ws.Range("D2:E51")
$gRange = $ws.Range("G2:G51")
$deRange.NumberFormat = "@"
$gRange.NumberFormat = "@"

$ws.Range("D2").Value = "305.06"
$ws.Range("E2").Value = "-1.47%"
$ws.Range("G2").Value = "18"
$ws.Range("D3").Value = "36.38"
$ws.Range("E3").Value = "-1.48%"
$ws.Range("G3").Value = "18"
$ws.Range("D4").Value = "5.030"
$ws.Range("E4").Value = "-0.19%"
$ws.Range("G4").Value = "18"
$ws.Range("D5").Value = "0.07904"
$ws.Range("E5").Value = "-0.15%"
$ws.Range("G5").Value = "18"
$ws.Range("D6").Value = "2.125"
$ws.Range("E6").Value = "-3.24%"
$ws.Range("G6").Value = "18"
$ws.Range("D7").Value = "7.964"
$ws.Range("E7").Value = "-0.89%"
$ws.Range("G7").Value = "18"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.138"
$ws.Range("E8").Value = "2.52%"
$ws.Range("G8").Value = "18"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9238"
$ws.Range("E9").Value = "-0.36%"
$ws.Range("G9").Value = "18"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.09663"
$ws.Range("E10").Value = "-2.79%"
$ws.Range("G10").Value = "18"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1860"
$ws.Range("E11").Value = "-1.64%"
$ws.Range("G11").Value = "18"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.08827"
$ws.Range("E12").Value = "1.53%"
$ws.Range("G12").Value = "18"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03602"
$ws.Range("E13").Value = "-0.06%"
$ws.Range("G13").Value = "18"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09932"
$ws.Range("E14").Value = "-0.28%"
$ws.Range("G14").Value = "18"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001432"
$ws.Range("E15").Value = "-3.73%"
$ws.Range("G15").Value = "18"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005627"
$ws.Range("E16").Value = "-1.04%"
$ws.Range("G16").Value = "18"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.472"
$ws.Range("E17").Value = "0.29%"
$ws.Range("G17").Value = "18"
$ws.Range("D18").Value = "2.653"
$ws.Range("E18").Value = "14.10%"
$ws.Range("G18").Value = "18"
$ws.Range("D19").Value = "0.3401"
$ws.Range("E19").Value = "-1.02%"
$ws.Range("G19").Value = "18"
$ws.Range("D20").Value = "0.1336"
$ws.Range("E20").Value = "0.76%"
$ws.Range("G20").Value = "18"
$ws.Range("D21").Value = "5.158"
$ws.Range("E21").Value = "4.43%"
$ws.Range("G21").Value = "18"
$ws.Range("D22").Value = "0.2246"
$ws.Range("E22").Value = "2.14%"
$ws.Range("G22").Value = "18"
$ws.Range("E23").Value = "-0.51%"
$ws.Range("G23").Value = "18"
$ws.Range("E24").Value = "-1.56%"
$ws.Range("G24").Value = "18"
$ws.Range("D25").Value = "0.004809"
$ws.Range("E25").Value = "-8.12%"
$ws.Range("G25").Value = "18"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").Value = "-7.13%"
$ws.Range("G26").Value = "18"
$ws.Range("D27").Value = "0.0004745"
$ws.Range("E27").Value = "74.76%"
$ws.Range("G27").Value = "18"
$ws.Range("G28").Value = "18"
$ws.Range("G29").Value = "18"
$ws.Range("G30").Value = "18"
$ws.Range("G31").Value = "18"
$ws.Range("G32").Value = "18"
$ws.Range("G33").Value = "18"
$ws.Range("G34").Value = "18"
$ws.Range("G35").Value = "18"
$ws.Range("G36").Value = "18"
$ws.Range("G37").Value = "18"
$ws.Range("G38").Value = "18"
$ws.Range("D39").Value = "0.01849"
$ws.Range("E39").Value = "1.06%"
$ws.Range("G39").Value = "18"
$ws.Range("D40").Value = "0.04898"
$ws.Range("E40").Value = "2.17%"
$ws.Range("G40").Value = "18"
$ws.Range("D41").Value = "0.007804"
$ws.Range("E41").Value = "-2.03%"
$ws.Range("G41").Value = "18"
$ws.Range("D42").Value = "0.1392"
$ws.Range("E42").Value = "-1.58%"
$ws.Range("G42").Value = "18"
$ws.Range("D43").Value = "0.007723"
$ws.Range("E43").Value = "2.49%"
$ws.Range("G43").Value = "18"
$ws.Range("D44").Value = "0.002207"
$ws.Range("E44").Value = "0.84%"
$ws.Range("G44").Value = "18"
$ws.Range("D45").Value = "0.01118"
$ws.Range("E45").Value = "10.38%"
$ws.Range("G45").Value = "18"
$ws.Range("D46").Value = "0.00006348"
$ws.Range("E46").Value = "2.08%"
$ws.Range("G46").Value = "18"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("G47").Value = "18"
$ws.Range("E48").Value = "0.26%"
$ws.Range("G48").Value = "18"
$ws.Range("D49").Value = "51.72"
$ws.Range("E49").Value = "43.29%"
$ws.Range("G49").Value = "18"
$ws.Range("D50").Value = "0.001899"
$ws.Range("E50").Value = "-29.35%"
$ws.Range("G50").Value = "18"
$ws.Range("D51").Value = "0.00002098"
$ws.Range("E51").Value = "-0.01%"
$ws.Range("G51").Value = "18"

$deRange.ClearFormats()
$gRange.ClearFormats()
